# "updated prompts and fixed edgecase in yes/no"
#  - B2/C2 (the _akhaliq row): Positive/Negative prompt text rewritten.
#  - Row 2 grows taller to fit the longer wrapped text.
#  - Selection cursor moves from B3 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Positive Prompt (B2) and Negative Prompt (C2) for _akhaliq ---

$b2Text = @'
- AI x Biology/Healthcare 
- Applications of multimodal transformer models and improvements to them
    - ex. LLM agent that uses vision to perform tasks on web browsers through screenshots (NOT normal transformer applications)
- Using LLMs to create new interfaces
    - ex. turning text into a graph or flow chart
- Startup launches
    - Only if they say how much they got funding (ex. Series A)
- Breakthroughs that are related to projects in my backlogs
    - AI text games combined with spaced repetition like Zork
    - Generating out-of-distribution ideas with LLMs
    - Music Generation
'@

$c2Text = @'
- Politics, culture, drama, news
    - Basically, things that I have no control over and won't help me
- Image creation models (ex. stable diffusion, midjourney, etc)
    - I am not working on craeting images right now
- Research that is not DIRECTLY useable to my current work and theoretical
    - ex. novel loss functions, pre-training hacks, large-scale training
    - Everything akhaliq posts is about AI so that is NOT a good reason in itself for reccomending a post. A lot of the machine learning research coming out isn't useable or relevant to the work I'm doing
- Time sinks that take me away from my work
    - I get nerd-sniped often, and it's one of the biggest things I need help with. If something doesn't fit the interested criteria, please shield me from getting nerd-sniped.
'@

# Trim the trailing newline that the here-string literal adds
$b2Text = $b2Text.TrimEnd("`r", "`n")
$c2Text = $c2Text.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $b2Text
$ws.Range("C2").Value = $c2Text

# Row 2 grew taller to fit the longer wrapped text
$ws.Rows.Item(2).RowHeight = 238

# --- Update the view / selection state ---
$ws.Range("C3").Select()

Write-Output "done"
